# Corrected info on firebase costs
#
# Slide 10 ("Cost") contains a paragraph describing Firebase pricing.
# The original text overstated the storage allowance and omitted the
# download bandwidth allowance, and described the $25 fee as yearly
# instead of monthly. Fix the wording without disturbing anything else
# on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$oldText = "otherwise `$25/year for 100k/instances for simultaneous connections and/or 50GB for storage."
$newText = "otherwise `$25/month for 100k/instances for simultaneous connections, 2.5GB for storage, and 20 GB downloads/month."

$result = $tr.Replace($oldText, $newText)
